# Add files via upload
# Applies the edits described by the commit diff:
#  - Sheet "Teste": header B1 text "Relevante/Irrelevante" -> "Relevancia"
#    (same shared string is reused, so this also seeds the new text for
#    the "Treinamento" header added below)
#  - Sheet "Treinamento": add new header cell B1 = "Relevancia"
#  - Sheet "Treinamento": flip several B-column flags from 0 -> 1 (and one 1 -> 0)
#  - Sheet "Teste": flip several B-column flags from 0 -> 1
#  - Sheet "Teste": set column A width (new <cols> block)
#  - Sheet selection/view bookkeeping: "Treinamento" becomes the active tab,
#    "Teste" keeps a fresh B1 selection but is no longer the active tab

$wb = $excel.ActiveWorkbook
$wsTreinamento = $wb.Worksheets.Item("Treinamento")
$wsTeste = $wb.Worksheets.Item("Teste")

# --- Header text / new header cell -----------------------------------------
# Re-text the existing shared string used by Teste!B1 ...
$wsTeste.Range("B1").Value = "Relevância"
# ... and add the same header to Treinamento!B1 (reuses the shared string)
$wsTreinamento.Range("B1").Value = "Relevância"

# --- Treinamento: B-column flag flips ---------------------------------------
$treinamentoRowsToOne = @(4,6,12,22,25,35,38,82,122,127,128,133,162,165,166,220,256,262,291,293,300)
foreach ($r in $treinamentoRowsToOne) {
    $wsTreinamento.Cells.Item($r, 2).Value = 1
}
$wsTreinamento.Cells.Item(287, 2).Value = 0

# --- Teste: B-column flag flips ---------------------------------------------
$testeRowsToOne = @(2,5,8,22,31,33,43,46,51,57,71,80,81,88,91,107,108,113,115,123,125,133,143,144,145,155,169,171,183,194,199)
foreach ($r in $testeRowsToOne) {
    $wsTeste.Cells.Item($r, 2).Value = 1
}

# --- Teste: widen column A ---------------------------------------------------
$wsTeste.Columns.Item(1).ColumnWidth = 105.5 - 5/6

# --- Selection / active-tab bookkeeping --------------------------------------
# Update Teste's selection first (this also makes it the active tab briefly) ...
$wsTeste.Range("B1").Select()
# ... then activate Treinamento last and select the full B column, so it ends
# up as the active tab, matching the final saved view state.
$wsTreinamento.Activate()
$wsTreinamento.Range("B1:B1048576").Select()
